# Update the table style applied to the six "member benefit" tables
# (old Table_0 custom style -> built-in table style) across the deck.

$p = $ppt.ActivePresentation

$oldStyleId = "{63E0C454-3585-45E7-966D-AA68C5263759}"
$newStyleId = "{E9145D0C-FFB7-4329-9AD6-17DDB7A5F2BF}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
